$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 794
$wsExhibit.Range("F6").Value = 11

# Sheet "全部类型" (All types) - duplicated rows for the same two events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 794
$wsAll.Range("F7").Value = 11
